# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (column G) values for rows 2-19 with the newly
# regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 2
    3  = 1
    4  = 5
    5  = 3
    6  = 6
    7  = 3
    8  = 8
    9  = 0
    10 = 3
    11 = 1
    12 = 1
    13 = 2
    14 = 7
    15 = 3
    16 = 7
    17 = 3
    18 = 6
    19 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
